# mailMap.xlsx — move the two-column mail map from A1:B2 down-and-right to
# B3:C4 (used to make room for a header/title block above and a leading
# column to the left), keep the hyperlink style/target on the moved
# "Théo Giraudet" e-mail cell, refresh the selection to the new block, and
# set an explicit (A4, portrait) page setup for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the data block (values + number formats + cell styles) -----------
# Cut/paste relocates the whole A1:B2 block so its top-left lands on B3,
# i.e. B3:C4 — this also keeps each moved cell's existing style (the
# hyperlink cell keeps its "Lien hypertexte" style without creating a
# duplicate style record).
$ws.Range("A1:B2").Cut($ws.Range("B3"))

# Cut leaves the old bottom-right corner (B2) behind as an empty, but still
# styled, cell — clear it so no stray formatting/cell entry remains.
$ws.Range("B2").Clear()

# --- Re-home the hyperlink --------------------------------------------------
# The hyperlink metadata itself doesn't follow Cut/Paste, so drop the old one
# (still anchored on B2) and recreate it on the relocated cell, C4, restoring
# the "Lien hypertexte" style the Add() call may disturb.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:theo.giraudet@etudiant.univ-rennes1.fr")
$ws.Range("C4").Style = "Lien hypertexte"

# --- Selection ---------------------------------------------------------------
$ws.Range("B3:C4").Select()

# --- Page setup: A4, portrait -----------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
